# Apply the account-list changes to the "Export" sheet.
# Operations are applied from the bottom of the sheet upward so that
# earlier (lower-numbered) row references stay valid as we go.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 125: delete 004948033 / GUILHERME / 69.7 (this account reappears
#    earlier in the sheet with an updated balance - see step 6).
$ws.Rows(125).Delete()

# 2) Row 17: delete 004550605 / REJANE / 869.47 (this account reappears
#    earlier in the sheet with an updated balance - see step 6).
$ws.Rows(17).Delete()

# 3) Insert a new row before row 12 (004457389 / RAFAEL) for the new
#    account 004207658 / ROBERTO / 937.62.
#    The account number is entered with a leading apostrophe (just like
#    typing it into Excel) so the leading zeros are kept and the cell
#    stays plain text instead of being coerced into a number.
$ws.Rows(12).Insert()
$ws.Range("A12").Formula = "'004207658"
$ws.Range("B12").Value = "ROBERTO"
$ws.Range("C12").Value = 937.62

# 4) Row 8: delete 004504449 / KELMA / 1000 (this account is re-inserted
#    one row earlier with an updated balance - see step 5).
$ws.Rows(8).Delete()

# 5) Insert a new row before row 7 (004488571 / CARLOS) for
#    004504449 / KELMA / 1024.15.
$ws.Rows(7).Insert()
$ws.Range("A7").Formula = "'004504449"
$ws.Range("B7").Value = "KELMA"
$ws.Range("C7").Value = 1024.15

# 6) Row 5 (005040864 / ANDRE / 12233.2) is replaced by two new rows:
#    004948033 / GUILHERME / 10311.28 and 004550605 / REJANE / 7445.66.
#    Insert a blank row first (pushing the old ANDRE row down to row 6),
#    then overwrite rows 5 and 6 with the new data.
$ws.Rows(5).Insert()
$ws.Range("A5").Formula = "'004948033"
$ws.Range("B5").Value = "GUILHERME"
$ws.Range("C5").Value = 10311.28
$ws.Range("A6").Formula = "'004550605"
$ws.Range("B6").Value = "REJANE"
$ws.Range("C6").Value = 7445.66

# 7) Row 2: delete 005547703 / SILVIA / 78861.27 entirely.
$ws.Rows(2).Delete()

# House-keeping: the engine leaves one stray blank row immediately above
# the trailing "Filtros aplicados..." footer row once row-shifting
# operations have run above. Drop that stray blank row so the footer
# sits directly after the last data row, exactly as in the original file.
$ws.Rows(224).Delete()
